$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 35
$ws.Range("B27").Value = 30.87
$ws.Range("C27").Value = 1.645
$ws.Range("D27").Value = 150.27

$ws.Range("A35").Value = 1
$ws.Range("B35").Value = 4.8
$ws.Range("C35").Value = 1.96

$ws.Range("A41").Value = 14
$ws.Range("B41").Value = 64
$ws.Range("C41").Value = 9
$ws.Range("D41").Value = 0.98
$ws.Range("E41").Value = 2.65

$ws.Range("A49").Value = 200
$ws.Range("B49").Value = 154
$ws.Range("E49").Value = 1.645

$ws.Range("A57").Value = 0.2
$ws.Range("B57").Value = 0.8
$ws.Range("C57").Value = 0.02
$ws.Range("D57").Value = 2.575

[void]$ws.Range("A49").Select()

